$d = $word.ActiveDocument

$d.Content.Find.Execute("419÷2=209, 1", $true, $false, $false, $false, $false, $true, 1, $false, "995÷3=331, 2", 2) | Out-Null
$d.Content.Find.Execute("662÷9=73, 5", $true, $false, $false, $false, $false, $true, 1, $false, "851÷4=212, 3", 2) | Out-Null
$d.Content.Find.Execute("188÷2=94, 0", $true, $false, $false, $false, $false, $true, 1, $false, "356÷2=178, 0", 2) | Out-Null
$d.Content.Find.Execute("984÷7=140, 4", $true, $false, $false, $false, $false, $true, 1, $false, "847÷8=105, 7", 2) | Out-Null
$d.Content.Find.Execute("710÷5=142, 0", $true, $false, $false, $false, $false, $true, 1, $false, "279÷9=31, 0", 2) | Out-Null
$d.Content.Find.Execute("923÷3=307, 2", $true, $false, $false, $false, $false, $true, 1, $false, "953÷5=190, 3", 2) | Out-Null
$d.Content.Find.Execute("706÷6=117, 4", $true, $false, $false, $false, $false, $true, 1, $false, "311÷8=38, 7", 2) | Out-Null
$d.Content.Find.Execute("892÷9=99, 1", $true, $false, $false, $false, $false, $true, 1, $false, "942÷2=471, 0", 2) | Out-Null
$d.Content.Find.Execute("482÷5=96, 2", $true, $false, $false, $false, $false, $true, 1, $false, "400÷5=80, 0", 2) | Out-Null
$d.Content.Find.Execute("549÷3=183, 0", $true, $false, $false, $false, $false, $true, 1, $false, "966÷8=120, 6", 2) | Out-Null
$d.Content.Find.Execute("306÷7=43, 5", $true, $false, $false, $false, $false, $true, 1, $false, "316÷7=45, 1", 2) | Out-Null
$d.Content.Find.Execute("448÷7=64, 0", $true, $false, $false, $false, $false, $true, 1, $false, "392÷9=43, 5", 2) | Out-Null
$d.Content.Find.Execute("401÷9=44, 5", $true, $false, $false, $false, $false, $true, 1, $false, "194÷6=32, 2", 2) | Out-Null
$d.Content.Find.Execute("421÷2=210, 1", $true, $false, $false, $false, $false, $true, 1, $false, "927÷5=185, 2", 2) | Out-Null
$d.Content.Find.Execute("427÷3=142, 1", $true, $false, $false, $false, $false, $true, 1, $false, "584÷3=194, 2", 2) | Out-Null
$d.Content.Find.Execute("209÷9=23, 2", $true, $false, $false, $false, $false, $true, 1, $false, "348÷5=69, 3", 2) | Out-Null
$d.Content.Find.Execute("207÷4=51, 3", $true, $false, $false, $false, $false, $true, 1, $false, "392÷2=196, 0", 2) | Out-Null
$d.Content.Find.Execute("181÷6=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "614÷2=307, 0", 2) | Out-Null
$d.Content.Find.Execute("673÷2=336, 1", $true, $false, $false, $false, $false, $true, 1, $false, "211÷6=35, 1", 2) | Out-Null
$d.Content.Find.Execute("425÷8=53, 1", $true, $false, $false, $false, $false, $true, 1, $false, "126÷4=31, 2", 2) | Out-Null
$d.Content.Find.Execute("780÷5=156, 0", $true, $false, $false, $false, $false, $true, 1, $false, "991÷3=330, 1", 2) | Out-Null
$d.Content.Find.Execute("325÷6=54, 1", $true, $false, $false, $false, $false, $true, 1, $false, "303÷9=33, 6", 2) | Out-Null
$d.Content.Find.Execute("443÷9=49, 2", $true, $false, $false, $false, $false, $true, 1, $false, "505÷6=84, 1", 2) | Out-Null
$d.Content.Find.Execute("653÷7=93, 2", $true, $false, $false, $false, $false, $true, 1, $false, "537÷5=107, 2", 2) | Out-Null
$d.Content.Find.Execute("811÷2=405, 1", $true, $false, $false, $false, $false, $true, 1, $false, "778÷8=97, 2", 2) | Out-Null
